# Refresh the crypto snapshot: Price (column D) and Volume(1h) change (column E)
# for each coin row, matching the latest scrape from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.075.76"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "2.306.65"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  +3.95%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "2.665.19"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "2.314.45"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "42.981.01"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.17%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("E30").Value = "  -10.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.40%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").Value = "2.008.16"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("D49").Value = "2.531.16"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.29%  "
$ws.Range("E51").Value = "  +1.62%  "
